$d = $word.ActiveDocument

# Locate the "GIS & Geospatial Analysis Consulting" paragraph under the
# Siege Analytics / PARTNER heading.
$anchor = $d.Paragraphs | Where-Object { $_.Range.Text -like "*GIS & Geospatial Analysis Consulting*" }

$lines = @(
    "• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels",
    "• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide",
    "• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis"
)

$current = $anchor
foreach ($line in $lines) {
    $current.Range.InsertParagraphAfter()
    $current = $current.Next()
    $current.Range.Text = $line
}

Write-Output "Inserted $($lines.Count) paragraphs after anchor."
